# feat: Set default value of `Show` to `true`.
#
# The "Show" column (F) in the data sheet should default to TRUE for the
# existing rows (F2:F13), displayed with a custom "TRUE/FALSE" boolean
# number format, and the current selection moved onto that range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$showRange = $ws.Range("F2:F13")

# Default the "Show" flag to TRUE for every existing data row.
$showRange.Value = $true

# Display booleans as the literal words TRUE / FALSE instead of 1/0.
$showRange.NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Leave the selection on the range that was just filled in (anchored at F3,
# matching what Excel does when a user drags/extends a selection down to F13).
$ws.Range("F3:F13").Select() | Out-Null
